$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "F" (Utilities/list) column. This shifts
# the old F,G,H columns to G,H,I and keeps their values/styles/validations intact.
$ws.Columns.Item(6).Insert()

# Populate the new F column with the shared "effective amount" formula.
$ws.Range("F1:F2").Formula = '=if(And(G1<>"",H1<>""),if(E1<>"",E1,D1),)'

# F1 takes on the shaded/right-aligned numeric look (reuse the existing
# shaded fill + Roboto font from the (now-shifted) old column H, I1, then
# layer on the 0.00 format).
$ws.Range("I1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").NumberFormat = "#,##0.00"
$ws.Range("F1").HorizontalAlignment = -4152

# F2 gets the same treatment but with the parenthesized-negative format.
$ws.Range("I2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").NumberFormat = "#,##0.00;(#,##0.00)"
$ws.Range("F2").HorizontalAlignment = -4152

$excel.CutCopyMode = 0
